$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K - strikeouts) values per row, replacing old Strike# based values
$gUpdates = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 2
    7 = 1
    8 = 0
    9 = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 3
    26 = 1
    27 = 2
    28 = 2
    29 = 0
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    40 = 1
    41 = 2
    42 = 1
    43 = 2
    44 = 1
    45 = 3
    46 = 0
    47 = 0
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 1
    54 = 1
    55 = 2
    56 = 2
    57 = 2
    58 = 1
    59 = 1
    60 = 2
    61 = 0
    62 = 1
    63 = 0
    64 = 2
    65 = 2
    66 = 0
    67 = 1
    68 = 1
    69 = 1
    70 = 0
    71 = 1
    72 = 0
    73 = 0
    74 = 2
}

foreach ($row in $gUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
}
